$d = $word.ActiveDocument

# 1. Update the title text: "Dataset and Data Analysis " -> "Visualization and Data Analysis."
$d.Content.Find.Execute("Dataset and Data Analysis ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Visualization and Data Analysis.", 2)

# 2. Move the _GoBack bookmark from its old location (an empty paragraph
#    near the end of the document) to right after the title paragraph's run.
$d.Bookmarks("_GoBack").Delete()

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $titleRange)
